$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data row (row 23) entirely so the used range / dimension shrinks back to A1:E5.
$ws.Rows.Item(23).Delete()

# Make sure the new row is stored as text (matches the original data, which was all text),
# even though some of the values look like numbers.
$ws.Range("A5:E5").NumberFormat = "@"

# Write the new data row into row 5.
$ws.Range("A5").Value = "sdadad 2313"
$ws.Range("B5").Value = "23132"
$ws.Range("C5").Value = "1231"
$ws.Range("D5").Value = "2312321"
$ws.Range("E5").Value = "stand"
